# tdf116486.docx: Open Sans -> DejaVu Sans for repeatable layout in test
#
# Six spots reference an "Open Sans*" family in word/styles.xml:
#   1. the document-wide rPrDefault (w:docDefaults/w:rPrDefault) - "Open Sans"
#   2. style "Normal"                       (Standard)                 - "Open Sans Light"
#   3. style "Fußbereich Standard"          (FubereichStandard)        - "Open Sans"
#   4. style "Strong"                       (Fett)                     - "Open Sans Semibold"
#   5. style "Querbalken 1. Ebene fett"     (Querbalken1Ebenefett)     - "Open Sans Semibold"
#   6. style "Seitenzahl Folgeseiten Zchn"  (SeitenzahlFolgeseitenZchn)- "Open Sans Light"
#
# All six become "DejaVu Sans" (ascii + hAnsi only - other rFonts attributes,
# eastAsiaTheme/cstheme/etc., are left untouched).

$d = $word.ActiveDocument

# 1) The five named styles that carry an explicit rFonts in their own rPr:
#    setting Font.Name rewrites just w:ascii/w:hAnsi on that style, leaving
#    every style that merely inherits "Open Sans" through the cascade alone.
$styleNames = @(
    "Normal",
    "Fußbereich Standard",
    "Strong",
    "Querbalken 1. Ebene fett",
    "Seitenzahl Folgeseiten Zchn"
)

foreach ($styleName in $styleNames) {
    $style = $d.Styles($styleName)
    $style.Font.Name = "DejaVu Sans"
}

# 2) The remaining spot is the package-wide default (w:docDefaults/w:rPrDefault),
#    which has no Style object of its own in the object model. Patch it via the
#    raw OOXML package, touching only that single element so nothing else moves.
$xml = $d.WordOpenXML

$pattern = '(<w:rPrDefault><w:rPr><w:rFonts w:ascii=")Open Sans("[^>]*w:hAnsi=")Open Sans("[^>]*/>)'
$replacement = '${1}DejaVu Sans${2}DejaVu Sans${3}'
$xml = [regex]::Replace($xml, $pattern, $replacement)

$d.WordOpenXML = $xml
